# Updates the addition/subtraction practice table.
#
# The edit removes the very first problem cell ("38-17=") from the table,
# shifting every following cell's answer "up" by one slot, and appends one
# brand-new problem ("3+10=" lands on the old last slot while a freshly
# generated value lands where the shift made room) -- net effect: every one
# of the 100 table cells (20 rows x 5 columns, read left-to-right/top-to-
# bottom) ends up holding a new value, taken from this list in document
# order.
$newValues = @(
  "59-50=","33-22=","68+30=","62-26=","85-37=","19+48=","34+23=","76-50=","79-62=","51-34=",
  "31+23=","22+42=","15-2=","84-57=","9+23=","10+18=","54-16=","43-35=","60+18=","80-7=",
  "75-23=","98-39=","87-38=","24-22=","0+49=","2+92=","47-34=","14+29=","44+35=","71-47=",
  "83-65=","39+30=","94-90=","29+23=","55+36=","57-53=","63+8=","82-9=","41+1=","37-27=",
  "1+62=","65-48=","9+8=","7+2=","58+36=","30+40=","20+67=","6-0=","68+24=","76-45=",
  "29-5=","49-18=","76-70=","89-47=","99-40=","23+41=","37-17=","68-8=","82+15=","22+61=",
  "55+19=","35-29=","62-17=","39+25=","5+48=","71+15=","27-18=","86+11=","96-62=","9+55=",
  "41-1=","71+8=","56-17=","59-25=","69-65=","13+21=","27-19=","87-87=","33-15=","86-16=",
  "50+32=","13+56=","73+26=","71-47=","81-44=","61+6=","74+5=","60-36=","35-11=","33+59=",
  "56-18=","71+25=","26+39=","92+1=","82-74=","18+35=","17+7=","67+14=","18+75=","3+10="
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}
